$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (Fecha = 2022-02-18, serial 44610) is inserted
# as the new row 56, pushing the former rows 56-59 down to 57-60.
# The new row duplicates every other field from the (former) row 56.
$ws.Rows.Item(56).Insert()

$ws.Cells.Item(56, 1).Value  = 1
$ws.Cells.Item(56, 2).Value  = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(56, 3).Value  = 'Arica y Parinacota'
$ws.Cells.Item(56, 4).Value  = 44610
$ws.Cells.Item(56, 5).Value  = 15
$ws.Cells.Item(56, 6).Value  = 100112021
$ws.Cells.Item(56, 7).Value  = 'Ají'
$ws.Cells.Item(56, 8).Value  = 'Inferno'
$ws.Cells.Item(56, 9).Value  = 'Primera'
$ws.Cells.Item(56, 10).Value = 120
$ws.Cells.Item(56, 11).Value = 14000
$ws.Cells.Item(56, 12).Value = 15000
$ws.Cells.Item(56, 13).Value = 14500
$ws.Cells.Item(56, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(56, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(56, 16).Value = 967
$ws.Cells.Item(56, 17).Value = 15
$ws.Cells.Item(56, 18).Value = 'Hortaliza'
